$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-25 20:58:47"

for ($row = 2; $row -le 72; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
